$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: phone number becomes a real number (drops leading zero)
$ws.Range("B7").Value = 753999382

# Row 8: new customer entry
$ws.Range("A8").Value = "chulanjana"
$ws.Range("B8").Value = 200031233443434
$ws.Range("C8").Value = "my phone is not working there is no signal at all"

# Row 9: new customer entry - phone number kept as text (leading zero preserved)
$ws.Range("A9").Value = "3323e32e342"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "0382250162"
$ws.Range("C9").Value = "my phone is not working"
